$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.595.30"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.918.92"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'245.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4832"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.83%  "
$ws.Range("D8").Value = "'0.2903"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "'0.06695"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").Value = "'107.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").Value = "'18.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "1.920.49"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "'0.07681"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "'5.290"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").Value = "'0.6687"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'279.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.75%  "
$ws.Range("D17").Value = "30.573.37"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007554"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.175.91"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'12.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").Value = "'5.500"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.27%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'6.401"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").Value = "'9.445"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "'164.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'20.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.99%  "
$ws.Range("D28").Value = "'2.123"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("D29").Value = "'0.1063"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "'1.405"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.04%  "
$ws.Range("D31").Value = "'4.153"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "'4.043"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").Value = "'0.05027"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("D34").Value = "'0.7325"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "'1.143"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'1.000"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.733"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.02033"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").Value = "'2.682"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'111.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("D41").Value = "'2.021"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").Value = "'0.4459"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("D43").Value = "'0.8743"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "'5.922"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'68.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'7.282"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.401"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'48.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.68%  "
$ws.Range("D50").Value = "'0.1255"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("D51").Value = "'35.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
